$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.144.73"
$ws.Range("E2").Value = "  +1.74%  "

$ws.Range("D3").Value = "3.917.11"
$ws.Range("E3").Value = "  +0.57%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "483.58"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +4.25%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "146.06"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.04%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.627"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.13%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.998"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.03%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.727"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.39%  "

$ws.Range("E10").Value = "  +2.16%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0000358"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +3.57%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "42.50"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.14%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "10.59"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.96%  "

$ws.Range("D14").Value = "4.543.62"
$ws.Range("E14").Value = "  +0.45%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "14.82"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.44%  "

$ws.Range("D16").Value = "3.891.32"
$ws.Range("E16").Value = "  -0.29%  "

$ws.Range("E17").Value = "  -0.24%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "19.81"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.77%  "

$ws.Range("E19").Value = "  -1.98%  "

$ws.Range("D20").Value = "68.309.39"
$ws.Range("E20").Value = "  +1.66%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "447.46"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.40%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "14.71"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.43%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.34"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.37%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "88.69"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.14%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "11.47"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +14.10%  "

$ws.Range("B26").Value = "PancakeSwap"
$ws.Range("C26").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.59"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.64%  "

$ws.Range("B27").Value = "RenderToken"
$ws.Range("C27").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.65"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +10.82%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "38.79"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.81%  "

$ws.Range("E29").Value = "  +2.65%  "

$ws.Range("B30").Value = "Bittensor"
$ws.Range("C30").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "694.27"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -6.01%  "

$ws.Range("B31").Value = "Cosmos"
$ws.Range("C31").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "13.43"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.53%  "

$ws.Range("B32").Value = "Hedera"
$ws.Range("C32").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.130"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.61%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.86"
$ws.Range("D33").Style = "Normal"

$ws.Range("E34").Value = "  +18.27%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "41.75"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -3.21%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "58.89"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.55%  "

$ws.Range("E37").Value = "  -4.82%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.63"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +5.35%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.998"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.02%  "

$ws.Range("E40").Value = "  +0.64%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.82"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +13.95%  "

$ws.Range("B42").Value = "ThetaToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.11"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.60%  "

$ws.Range("B43").Value = "TheGraph"
$ws.Range("C43").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.363"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +8.44%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.98"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +6.25%  "

$ws.Range("E45").Value = "  +0.98%  "

$ws.Range("E46").Value = "  -0.14%  "

$ws.Range("E47").Value = "  +0.34%  "

$ws.Range("E48").Value = "  -1.58%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "145.82"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.75%  "

$ws.Range("E50").Value = "  -1.17%  "

$ws.Range("E51").Value = "  -1.12%  "

